$d = $word.ActiveDocument

# The document currently ends right after the last </w:tbl> (no trailing
# paragraph). $d.Paragraphs.Last is the zero-length "end of story" range
# that sits right after that last table; replacing its contents with raw
# OOXML appends the new content immediately after the existing last table.
$endRange = $d.Paragraphs.Last.Range

$newTableXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblW w:w="2500" w:type="pct"/>
    <w:jc w:val="center"/>
    <w:tblBorders>
      <w:top w:val="single" w:color="auto" w:sz="8"/>
      <w:left w:val="single" w:color="auto" w:sz="8"/>
      <w:bottom w:val="single" w:color="auto" w:sz="8"/>
      <w:right w:val="single" w:color="auto" w:sz="8"/>
      <w:insideH w:val="single" w:color="auto" w:sz="8"/>
      <w:insideV w:val="single" w:color="auto" w:sz="8"/>
    </w:tblBorders>
  </w:tblPr>
  <w:tblGrid/>
  <w:tr>
    <w:trPr>
      <w:trHeight w:val="2880" w:hRule="exact"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:vAlign w:val="center"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="center"/>
        </w:pPr>
        <w:r>
          <w:t>hello world</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
"@

$endRange.InsertXML($newTableXml)
